# Apply updated crypto price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.303.09"
$ws.Range("E2").Value = "  -2.42%  "
$ws.Range("D3").Value = "2.991.83"
$ws.Range("E3").Value = "  -3.07%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").Value = "'584.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("D6").Value = "'145.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.68%  "
$ws.Range("D8").Value = "'0.521"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.18%  "
$ws.Range("D9").Value = "2.992.48"
$ws.Range("E9").Value = "  -2.97%  "
$ws.Range("E10").Value = "  -6.51%  "
$ws.Range("E11").Value = "  -4.26%  "
$ws.Range("D12").Value = "'0.441"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.34%  "
$ws.Range("D13").Value = "'0.0000227"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.85%  "
$ws.Range("D14").Value = "'34.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.36%  "
$ws.Range("D15").Value = "'0.123"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.07%  "
$ws.Range("D16").Value = "3.477.74"
$ws.Range("E16").Value = "  -3.34%  "
$ws.Range("D17").Value = "62.243.29"
$ws.Range("E17").Value = "  -2.43%  "
$ws.Range("D18").Value = "'6.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.38%  "
$ws.Range("D19").Value = "2.988.88"
$ws.Range("E19").Value = "  -3.48%  "
$ws.Range("D20").Value = "'456.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.74%  "
$ws.Range("D21").Value = "'13.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.55%  "
$ws.Range("D22").Value = "'0.676"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.09%  "
$ws.Range("D23").Value = "'7.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.93%  "
$ws.Range("D24").Value = "'80.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.89%  "
$ws.Range("D25").Value = "'2.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.19%  "
$ws.Range("D26").Value = "'12.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.30%  "
$ws.Range("D27").Value = "'10.07"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.57%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'7.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.52%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'2.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.38%  "
$ws.Range("E32").Value = "  -5.70%  "
$ws.Range("D33").Value = "'26.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").Value = "'0.108"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.24%  "
$ws.Range("E35").Value = "  -5.04%  "
$ws.Range("D36").Value = "0.0₃0781"
$ws.Range("E36").Value = "  -7.05%  "
$ws.Range("E37").Value = "  -5.34%  "
$ws.Range("D38").Value = "'2.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.58%  "
$ws.Range("D39").Value = "'50.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("D40").Value = "'8.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.65%  "
$ws.Range("E41").Value = "  -11.39%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.113"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'388.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -12.00%  "
$ws.Range("D44").Value = "'0.272"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.99%  "
$ws.Range("D45").Value = "2.757.59"
$ws.Range("E45").Value = "  -2.73%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").Value = "'0.0349"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.03%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").Value = "'38.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.24%  "
$ws.Range("D48").Value = "'127.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.83%  "
$ws.Range("E50").Value = "  -2.10%  "
$ws.Range("D51").Value = "'23.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.03%  "
